$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E10').Value = 'satisf pel reconhec autocobranç par correspond expect cans somat'
$ws.Range('E11').Value = 'satisf pel reconhec autocobranç par correspond expect cans somat'
$ws.Range('E12').Value = 'satisf pel reconhec autocobranç par correspond expect cans somat'
$ws.Range('A13').Value = 'coleg telef ped ajud últ hor visit famili trabalh faculdad par'
$ws.Range('A14').Value = 'coleg telef ped ajud últ hor visit famili trabalh faculdad par'
$ws.Range('A17').Value = 'restr aliment cas so disponibil discuss coleg cl'
$ws.Range('A18').Value = 'restr aliment cas so disponibil discuss coleg cl'
$ws.Range('A19').Value = 'restr aliment cas so disponibil discuss coleg cl'
$ws.Range('A20').Value = 'restr aliment cas so disponibil discuss coleg cl'
$ws.Range('C24').Value = 'contracontrol control avers pal'
$ws.Range('A25').Value = 'priv vint hor est so cas autorregr vou alg calór grand quant par consegu vomit facil discuss pai'
$ws.Range('A26').Value = 'priv vint hor est so cas autorregr vou alg calór grand quant par consegu vomit facil discuss pai'
$ws.Range('A27').Value = 'priv vint hor est so cas autorregr vou alg calór grand quant par consegu vomit facil discuss pai'
$ws.Range('A28').Value = 'priv vint hor est so cas autorregr vou alg calór grand quant par consegu vomit facil discuss pai'
$ws.Range('B29').Value = 'induz vômit us lax us diuré restring aliment exercíci físic'
$ws.Range('B30').Value = 'induz vômit us lax us diuré restring aliment exercíci físic'
$ws.Range('B31').Value = 'induz vômit us lax us diuré restring aliment exercíci físic'
$ws.Range('E31').Value = 'autoconfianç par invest relacion amor'
$ws.Range('B32').Value = 'induz vômit us lax us diuré restring aliment exercíci físic'
$ws.Range('E32').Value = 'autoconfianç par invest relacion amor'
$ws.Range('B33').Value = 'induz vômit us lax us diuré restring aliment exercíci físic'
$ws.Range('E33').Value = 'autoconfianç par invest relacion amor'
$ws.Range('B34').Value = 'induz vômit us lax us diuré restring aliment exercíci físic'
$ws.Range('E34').Value = 'autoconfianç par invest relacion amor'
$ws.Range('C42').Value = 'contracontrol control avers pai'
$ws.Range('A60').Value = 'discuss delegac comum pai vítim estrup tent atir acus dur discuss revolv cai próx vítim golp sit relacion afet amor vítim peg revólv apont par cabeç'
$ws.Range('A61').Value = 'discuss delegac comum pai vítim estrup tent atir acus dur discuss revolv cai próx vítim golp sit relacion afet amor vítim peg revólv apont par cabeç'
$ws.Range('C61').Value = 'verôn presenc suicídi mulh trabalh'
$ws.Range('A62').Value = 'discuss delegac comum pai vítim estrup tent atir acus dur discuss revolv cai próx vítim golp sit relacion afet amor vítim peg revólv apont par cabeç'
$ws.Range('C66').Value = 'chef diss pod investig cas porqu presenci suicídi escrivã detetiv'
$ws.Range('A72').Value = 'verôn descobr corrupç polic regr polic ajud pesso'
$ws.Range('A73').Value = 'verôn descobr corrupç polic regr polic ajud pesso'
$ws.Range('A79').Value = 'inter soc convit par situ soc demand tom decis problem resolv nov oportun trabalh'
$ws.Range('B79').Value = 'respost passiv tímid p ex fal pouc pens adi decis sab med julg decis encerr assunt respond med ansiedad'
$ws.Range('A80').Value = 'inter soc convit par situ soc demand tom decis problem resolv nov oportun trabalh'
$ws.Range('B80').Value = 'respost passiv tímid p ex fal pouc pens adi decis sab med julg decis encerr assunt respond med ansiedad'
$ws.Range('A81').Value = 'inter soc convit par situ soc demand tom decis problem resolv nov oportun trabalh'
$ws.Range('B81').Value = 'respost passiv tímid p ex fal pouc pens adi decis sab med julg decis encerr assunt respond med ansiedad'
$ws.Range('A82').Value = 'inter soc convit par situ soc demand tom decis problem resolv nov oportun trabalh'
$ws.Range('B82').Value = 'respost passiv tímid p ex fal pouc pens adi decis sab med julg decis encerr assunt respond med ansiedad'
$ws.Range('B83').Value = 'assum respons pel resoluç tod quest relacion filh desmarc compromiss trabalh desmarc ativ laz'
$ws.Range('B84').Value = 'assum respons pel resoluç tod quest relacion filh desmarc compromiss trabalh desmarc ativ laz'
$ws.Range('B85').Value = 'assum respons pel resoluç tod quest relacion filh desmarc compromiss trabalh desmarc ativ laz'
$ws.Range('B86').Value = 'assum respons pel resoluç tod quest relacion filh desmarc compromiss trabalh desmarc ativ laz'
$ws.Range('B87').Value = 'assum respons pel resoluç tod quest relacion filh desmarc compromiss trabalh desmarc ativ laz'
$ws.Range('B88').Value = 'assum respons pel resoluç tod quest relacion filh desmarc compromiss trabalh desmarc ativ laz'
$ws.Range('C97').Value = 'sobrecarg pesso aproveit del'
$ws.Range('A100').Value = 'convit namor par viaj convit amig coleg trabalh par sair propost trabalh problem saúd demand relat cuid filh'
$ws.Range('B100').Value = 'aceit convit par sair ativ difer aceit nov propost trabalh ped ajud cuid filh divid respons cuid saúd aliment médic adequ pratic ativ físic'
$ws.Range('A101').Value = 'convit namor par viaj convit amig coleg trabalh par sair propost trabalh problem saúd demand relat cuid filh'
$ws.Range('B101').Value = 'aceit convit par sair ativ difer aceit nov propost trabalh ped ajud cuid filh divid respons cuid saúd aliment médic adequ pratic ativ físic'
$ws.Range('A102').Value = 'convit namor par viaj convit amig coleg trabalh par sair propost trabalh problem saúd demand relat cuid filh'
$ws.Range('B102').Value = 'aceit convit par sair ativ difer aceit nov propost trabalh ped ajud cuid filh divid respons cuid saúd aliment médic adequ pratic ativ físic'
$ws.Range('A103').Value = 'convit namor par viaj convit amig coleg trabalh par sair propost trabalh problem saúd demand relat cuid filh'
$ws.Range('B103').Value = 'aceit convit par sair ativ difer aceit nov propost trabalh ped ajud cuid filh divid respons cuid saúd aliment médic adequ pratic ativ físic'
$ws.Range('B104').Value = 'assert argument express opin sent recus ped form diret'
$ws.Range('B105').Value = 'assert argument express opin sent recus ped form diret'
$ws.Range('B106').Value = 'assert argument express opin sent recus ped form diret'
$ws.Range('B107').Value = 'assert argument express opin sent recus ped form diret'
$ws.Range('B108').Value = 'assert argument express opin sent recus ped form diret'
$ws.Range('B109').Value = 'assert argument express opin sent recus ped form diret'
$ws.Range('B110').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B111').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('C111').Value = 'risc pesso recus ped'
$ws.Range('B112').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B113').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B114').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B115').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B116').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B117').Value = 'busc soluç solic ajud ped diret envolv nov oportun trabalh tom decis express sent necess man cl'
$ws.Range('B121').Value = 'hiperlal'
$ws.Range('B122').Value = 'hiperlal'
$ws.Range('A128').Value = 'regr regr pai soub alg sobr voc mat pesso famíl unid brig autorregr pesso tod dev gost mim sd solicit ajud amig famili'
$ws.Range('A129').Value = 'regr regr pai soub alg sobr voc mat pesso famíl unid brig autorregr pesso tod dev gost mim sd solicit ajud amig famili'
$ws.Range('A130').Value = 'regr regr pai soub alg sobr voc mat pesso famíl unid brig autorregr pesso tod dev gost mim sd solicit ajud amig famili'
$ws.Range('A131').Value = 'regr regr pai soub alg sobr voc mat pesso famíl unid brig autorregr pesso tod dev gost mim sd solicit ajud amig famili'
$ws.Range('B149').Value = 'ativ norm dentr aceit pel pal evit pai fing vê escut respost respond ansiedad'
$ws.Range('B150').Value = 'ativ norm dentr aceit pel pal evit pai fing vê escut respost respond ansiedad'
$ws.Range('B156').Value = 'aument voz xing discut car raiv consegu pens cl respost respond transpir sens est estress'
$ws.Range('B157').Value = 'aument voz xing discut car raiv consegu pens cl respost respond transpir sens est estress'
$ws.Range('B158').Value = 'aument voz xing discut car raiv consegu pens cl respost respond transpir sens est estress'
$ws.Range('B183').Value = 'adi comport procrastin vésp tent lembr acontec esquec pouc'
$ws.Range('B184').Value = 'adi comport procrastin vésp tent lembr acontec esquec pouc'
